$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 8, pushing existing rows 8-69 down to 9-70.
$ws.Rows.Item(8).Insert()

# The sheet keeps an explicit (empty) cell in every column A:Y for every
# data row, so stamp the same empty-string pattern into the freshly
# inserted row before filling in the real values.
for ($col = 1; $col -le 25; $col++) {
    $ws.Cells.Item(8, $col).Value = ""
}

# Populate the newly inserted row 8 with the new September transaction.
$ws.Cells.Item(8, 18).Value = "teknothermindustries anyone axis"
$ws.Cells.Item(8, 19).Value = "2024-09-05 16:18:13"
